# Apply the coin-price / volume refresh described by the commit
# "Updated cryptos list on Wed May 31 11:47:04 UTC 2023 with GitHub Actions"
#
# Each cell is forced to the Text number format before the write (and reverted
# to the workbook's Normal style afterwards) so that values such as "5.315" or
# "27.119.16" are stored as literal text -- matching the source data -- instead
# of being auto-coerced into numbers by Excel's smart-entry parser.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "27.119.16"
Set-TextValue "E2" "  -2.91%  "
# Row 3
Set-TextValue "D3" "1.867.62"
Set-TextValue "E3" "  -2.25%  "
# Row 4
Set-TextValue "E4" "  +0.19%  "
# Row 5
Set-TextValue "D5" "306.83"
Set-TextValue "E5" "  -2.00%  "
# Row 6
Set-TextValue "E6" "  +0.13%  "
# Row 7
Set-TextValue "D7" "0.5131"
Set-TextValue "E7" "  +2.40%  "
# Row 8
Set-TextValue "D8" "0.3743"
Set-TextValue "E8" "  -1.95%  "
# Row 9
Set-TextValue "D9" "0.07136"
Set-TextValue "E9" "  -2.36%  "
# Row 10
Set-TextValue "D10" "0.8885"
Set-TextValue "E10" "  -2.65%  "
# Row 11
Set-TextValue "D11" "20.66"
Set-TextValue "E11" "  -2.80%  "
# Row 12
Set-TextValue "D12" "0.07528"
Set-TextValue "E12" "  -1.85%  "
# Row 13
Set-TextValue "D13" "1.834.12"
Set-TextValue "E13" "  -5.65%  "
# Row 14
Set-TextValue "D14" "5.315"
Set-TextValue "E14" "  -3.09%  "
# Row 15
Set-TextValue "D15" "89.24"
Set-TextValue "E15" "  -3.86%  "
# Row 16
Set-TextValue "E16" "  +0.20%  "
# Row 17
Set-TextValue "D17" "0.000008478"
Set-TextValue "E17" "  -3.05%  "
# Row 18
Set-TextValue "E18" "  -3.70%  "
# Row 19
Set-TextValue "D19" "0.9999"
Set-TextValue "E19" "  +0.18%  "
# Row 20
Set-TextValue "D20" "27.170.84"
Set-TextValue "E20" "  -2.85%  "
# Row 21
Set-TextValue "D21" "5.051"
Set-TextValue "E21" "  -2.60%  "
# Row 22
Set-TextValue "D22" "2.082.50"
Set-TextValue "E22" "  -3.56%  "
# Row 23
Set-TextValue "E23" "  -2.67%  "
# Row 24
Set-TextValue "D24" "6.479"
Set-TextValue "E24" "  -2.09%  "
# Row 25
Set-TextValue "D25" "149.87"
Set-TextValue "E25" "  -2.16%  "
# Row 26
Set-TextValue "E26" "  +0.06%  "
# Row 28
Set-TextValue "D28" "2.100"
Set-TextValue "E28" "  -4.90%  "
# Row 29
Set-TextValue "D29" "112.82"
Set-TextValue "E29" "  -2.25%  "
# Row 30
Set-TextValue "D30" "4.761"
Set-TextValue "E30" "  -3.46%  "
# Row 31
Set-TextValue "D31" "4.677"
Set-TextValue "E31" "  -3.56%  "
# Row 32
Set-TextValue "D32" "0.09036"
Set-TextValue "E32" "  +0.13%  "
# Row 33
Set-TextValue "D33" "0.05137"
Set-TextValue "E33" "  -2.84%  "
# Row 34
Set-TextValue "D34" "3.095"
Set-TextValue "E34" "  -3.56%  "
# Row 35
Set-TextValue "E35" "  -6.18%  "
# Row 36
Set-TextValue "D36" "0.7369"
Set-TextValue "E36" "  -5.77%  "
# Row 37
Set-TextValue "D37" "0.02043"
Set-TextValue "E37" "  -2.15%  "
# Row 38
Set-TextValue "D38" "2.504"
Set-TextValue "E38" "  -4.07%  "
# Row 39
Set-TextValue "D39" "3.044"
Set-TextValue "E39" "  -0.83%  "
# Row 40
Set-TextValue "E40" "  -1.83%  "
# Row 41
Set-TextValue "D41" "0.5298"
Set-TextValue "E41" "  -4.63%  "
# Row 42
Set-TextValue "D42" "6.601"
Set-TextValue "E42" "  -4.18%  "
# Row 43
Set-TextValue "D43" "116.64"
Set-TextValue "E43" "  +2.57%  "
# Row 44
Set-TextValue "D44" "8.333"
Set-TextValue "E44" "  -2.47%  "
# Row 45
Set-TextValue "D45" "0.1474"
Set-TextValue "E45" "  -2.91%  "
# Row 46
Set-TextValue "B46" "PaxDollar"
Set-TextValue "C46" "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue "D46" "0.9993"
Set-TextValue "E46" "  +0.15%  "
# Row 47
Set-TextValue "B47" "Decentraland"
Set-TextValue "C47" "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue "D47" "0.4626"
Set-TextValue "E47" "  -4.27%  "
# Row 48
Set-TextValue "D48" "10.03"
Set-TextValue "E48" "  -5.95%  "
# Row 49
Set-TextValue "D49" "1.571"
Set-TextValue "E49" "  -4.27%  "
# Row 50
Set-TextValue "D50" "64.46"
Set-TextValue "E50" "  -4.74%  "
# Row 51
Set-TextValue "D51" "36.54"
Set-TextValue "E51" "  -1.63%  "
